# Apply updated cryptocurrency price/volume data to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.873.17'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '2.502.50'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.524'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.555'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.95'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.125'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.78'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').Value = '2.902.20'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '2.508.15'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.856'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '47.819.74'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('E19').Value = '  +2.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.64'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +14.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '247.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.94'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.141'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0793'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.49'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.76%  '
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0300'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').Value = '2.005.88'
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.23%  '
